$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$filesQuery = @'
MATCH (f:file)-->(parent)
MATCH (f)-[:file_of_sample]->(samp)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)-[:study_subject_of_study]->(s)
MATCH (s)-[:study_of_program]->(p)
MATCH (d)-[:diagnosis_of_study_subject]->(ss)
MATCH (tp)-[:tp_of_diagnosis]->(d)
WHERE  f.file_format IN ['gz','tbi']
WITH
        f, parent,p, ss, d,tp, s, samp,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent,p, ss, d,tp, s, samp,
        f.file_size /(1024^i) AS value,
        10^precision AS factor,
        units[i] as unit
WITH
        f, parent,p, ss, d,tp, s, samp, unit,
        round(factor * value)/factor AS size
RETURN Distinct
    f.file_name AS `File Name`,
    head(labels(samp)) AS `Association`,
    f.file_description AS `Description`,
    f.file_format AS `File Format`,
     CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
    p.program_acronym AS `Program Code`,
    s.study_acronym AS `Arm`,
    ss.study_subject_id AS `Case ID`,
    samp.sample_id AS `Sample ID`
     order By f.file_name ASC LIMIT 100
'@

$samplesQuery = @'
MATCH (ss:study_subject)
WITH COLLECT(ss.study_subject_id) AS all_subjects
MATCH (samp:sample)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (samp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
WHERE  f.file_format IN ['gz','tbi']
WITH
    distinct lp,
    toInteger(split(ss.study_subject_id,'-')[2]) AS subject_id_num,
    collect(distinct f.file_id) AS files,
    samp, ss, s, p, all_subjects
RETURN
 samp.sample_id AS `Sample ID`,
            ss.study_subject_id AS `Case ID`,
            p.program_acronym AS `Program Code`,
            s.study_acronym AS `Arm`,
            ss.disease_subtype AS `Diagnosis`,
            samp.tissue_type AS `Tissue Type`,
            samp.composition AS `Tissue Composition`,
            samp.sample_anatomic_site AS `Sample Anatomic Site`,
            samp.method_of_sample_procurement AS `Sample Procurement Method`
 order By samp.sample_id ASC LIMIT 100
'@

$casesQuery = @'
MATCH (ss:study_subject)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)<-[:sample_of_study_subject]-(samp:sample)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)

WITH DISTINCT ss, samp, collect(DISTINCT samp.sample_id) AS samples, collect(DISTINCT lp.laboratory_procedure_id) AS lab_procedures, collect(DISTINCT f) AS files
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
MATCH (ss)<-[:demographic_of_study_subject]-(demo)
MATCH (f:file)
MATCH (f)-[:file_of_sample]->(samp)
WHERE  f.file_format IN ['gz','tbi']
return DISTINCT ss.study_subject_id as `Case ID`,
   p.program_acronym as `Program Code`,
    p.program_id as Program_ID,
   s.study_acronym as `Arm`,
   ss.disease_subtype as `Diagnosis`,
   sf.grouped_recurrence_score AS `Recurrence Score`,
   d.tumor_size_group AS `tumor_size`,
   d.er_status AS `ER Status`,
   d.pr_status AS `PR Status`,
   coalesce(CASE demo.age_at_index % 1 WHEN 0 THEN apoc.convert.toInteger(demo.age_at_index) ELSE demo.age_at_index END, '') AS `Age (years)`,
demo.survival_time AS `Survival (days)`
 order By ss.study_subject_id ASC LIMIT 100
'@

# Order matters: the shared-string table assigns new indices in the order
# cells are written, and row 4 (Files) must land at the lowest index (11),
# row 3 (Samples) next (12), row 2 (Cases) last (13) to match the target file.
$ws.Range("B4").Value = $filesQuery
$ws.Range("B3").Value = $samplesQuery
$ws.Range("B2").Value = $casesQuery

# The wrapped-text queries now span one extra line each (except Files, whose
# text only grew on its existing last line) so rows 2 and 3 grow taller.
$ws.Rows.Item(2).RowHeight = 403.2
$ws.Rows.Item(3).RowHeight = 360

$ws.Range("B2").Select() | Out-Null
